# daily auto push: 2026-02-27 09:56 UTC
#
# The source sheet ("日付/曜日/時刻/ランキング" log) gets one new sample
# appended for "today" (2026/02/27) and inserted in date order. That pushes
# the previously-last block of rows (2026/12/29 .. 2027/01/05) down by one
# row, growing the used range from A1:D911 to A1:D912.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 870

# Make room for the new sample: shifts rows 870-911 down to 871-912.
$ws.Rows.Item($newRow).Insert()

# Column A holds dates as plain text (e.g. "2026/12/29"), not real Excel
# dates, throughout this sheet. A bare .Value assignment of a
# yyyy/mm/dd-looking string gets auto-parsed into a date serial by Excel,
# so prefix with a quote to force literal text, matching every other row,
# then drop the formatting that the quote-prefix leaves behind.
$ws.Cells.Item($newRow, 1).Value = "'2026/02/27"
$ws.Cells.Item($newRow, 2).Value = "金"
$ws.Cells.Item($newRow, 3).Value = 16
$ws.Cells.Item($newRow, 4).Value = 44

$ws.Range("A" + $newRow + ":B" + $newRow).ClearFormats()
